$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 91
$ws1.Range("F3").Value = 4040
$ws1.Range("F11").Value = 78
$ws1.Range("F13").Value = 1508
$ws1.Range("F14").Value = 268
$ws1.Range("F15").Value = 2861

# Sheet "全部类型" (All categories)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 91
$ws4.Range("F3").Value = 4040
$ws4.Range("F12").Value = 78
$ws4.Range("F16").Value = 1508
$ws4.Range("F17").Value = 268
$ws4.Range("F18").Value = 2861
